# Updated cryptos list on Sat Oct  7 08:36:53 UTC 2023 with GitHub Actions
#
# Refreshes the Price (column D) and Volume(1h) (column E) figures for the
# crypto table on the active sheet, and fixes the rank ordering of rows
# 43/44 (RenderToken now ranks ahead of Aave).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price/Volume columns are plain text in the source data (e.g. prices
# like "1.640.27" use '.' as a thousands separator). Force the cells back
# to text formatting first so Excel doesn't auto-coerce numeric-looking
# values (e.g. "213.42", "0.920") into actual numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "27.930.43"
$ws.Range("E2").Value = "  +1.14%  "

$ws.Range("D3").Value = "1.640.82"
$ws.Range("E3").Value = "  +0.40%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").Value = "213.42"
$ws.Range("E5").Value = "  +0.49%  "

$ws.Range("D6").Value = "0.524"
$ws.Range("E6").Value = "  +0.46%  "

$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.04%  "

$ws.Range("D8").Value = "23.71"
$ws.Range("E8").Value = "  +1.55%  "

$ws.Range("D9").Value = "0.262"
$ws.Range("E9").Value = "  -0.78%  "

$ws.Range("D10").Value = "0.0615"
$ws.Range("E10").Value = "  +0.50%  "

$ws.Range("D11").Value = "0.0877"
$ws.Range("E11").Value = "  +0.71%  "

$ws.Range("D12").Value = "1.872.53"
$ws.Range("E12").Value = "  +0.28%  "

$ws.Range("D13").Value = "1.640.56"
$ws.Range("E13").Value = "  +0.30%  "

$ws.Range("D14").Value = "4.10"
$ws.Range("E14").Value = "  +1.20%  "

$ws.Range("E15").Value = "  +3.79%  "

$ws.Range("D16").Value = "65.98"
$ws.Range("E16").Value = "  +1.06%  "

$ws.Range("D17").Value = "27.916.71"
$ws.Range("E17").Value = "  +1.10%  "

$ws.Range("D18").Value = "231.31"
$ws.Range("E18").Value = "  +0.02%  "

$ws.Range("E19").Value = "  +0.59%  "

$ws.Range("D20").Value = "7.61"
$ws.Range("E20").Value = "  +0.39%  "

$ws.Range("E21").Value = "  +0.07%  "

$ws.Range("D22").Value = "10.77"
$ws.Range("E22").Value = "  +1.34%  "

$ws.Range("E23").Value = "  +0.17%  "

$ws.Range("E24").Value = "  -3.06%  "

$ws.Range("D25").Value = "151.51"
$ws.Range("E25").Value = "  +1.21%  "

$ws.Range("D26").Value = "6.95"
$ws.Range("E26").Value = "  +0.81%  "

$ws.Range("D27").Value = "15.69"
$ws.Range("E27").Value = "  +0.83%  "

$ws.Range("E28").Value = "  -0.07%  "

$ws.Range("E29").Value = "  -0.01%  "

$ws.Range("E30").Value = "  +0.78%  "

$ws.Range("E31").Value = "  +0.06%  "

$ws.Range("E32").Value = "  +1.64%  "

$ws.Range("E33").Value = "  +0.59%  "

$ws.Range("D34").Value = "1.397.73"
$ws.Range("E34").Value = "  -5.45%  "

$ws.Range("E35").Value = "  +1.53%  "

$ws.Range("E36").Value = "  +0.58%  "

$ws.Range("E37").Value = "  +0.93%  "

$ws.Range("D38").Value = "0.0168"
$ws.Range("E38").Value = "  +0.15%  "

$ws.Range("D39").Value = "0.920"
$ws.Range("E39").Value = "  -2.27%  "

$ws.Range("D40").Value = "0.556"
$ws.Range("E40").Value = "  -0.67%  "

$ws.Range("E41").Value = "  -0.65%  "

$ws.Range("E42").Value = "  +0.01%  "

# Row 43/44: RenderToken and Aave swap rank positions.
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").Value = "1.84"
$ws.Range("E43").Value = "  +5.16%  "

$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").Value = "66.27"
$ws.Range("E44").Value = "  -2.66%  "

$ws.Range("D45").Value = "5.44"
$ws.Range("E45").Value = "  +1.48%  "

$ws.Range("D46").Value = "2.21"
$ws.Range("E46").Value = "  +0.08%  "

$ws.Range("D47").Value = "1.781.30"
$ws.Range("E47").Value = "  +0.52%  "

$ws.Range("D48").Value = "87.99"
$ws.Range("E48").Value = "  +0.48%  "

$ws.Range("E49").Value = "  +0.88%  "

$ws.Range("E50").Value = "  +0.35%  "

$ws.Range("D51").Value = "7.64"
$ws.Range("E51").Value = "  -1.29%  "
